$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update control-signal table values (H6, I6, H7, I7, H8)
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 1
$ws.Range("H8").Value = 1

# Update the active cell selection on the sheet
$ws.Activate()
$ws.Range("J7").Select()
